# Applies the cryptos-list refresh: updates Price (column D) and
# Volume(1h) (column E) text values for the affected rows.
#
# The source cells are stored as text (inline strings), and several of the
# new Price values look like plain numbers (e.g. "1.006", "21.91"). Assigning
# those directly to Range.Value would make Excel auto-convert them to numeric
# values (losing trailing zeros / introducing floating-point artifacts), so
# such values are prefixed with a leading apostrophe, matching how Excel's
# UI forces text-entry for number-looking strings; the apostrophe itself is
# not stored in the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.238.09'
$ws.Range("E2").Value = '  -7.78%  '
$ws.Range("D3").Value = '1.676.10'
$ws.Range("E3").Value = '  -5.38%  '
$ws.Range("E4").Value = '  +0.29%  '
$ws.Range("E5").Value = '  -4.57%  '
$ws.Range("D6").Value = '''0.5125'
$ws.Range("E6").Value = '  -12.52%  '
$ws.Range("D7").Value = '''1.006'
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").Value = '''0.2661'
$ws.Range("E8").Value = '  -2.71%  '
$ws.Range("D9").Value = '''21.91'
$ws.Range("E9").Value = '  -5.29%  '
$ws.Range("D10").Value = '''0.06340'
$ws.Range("E10").Value = '  -4.98%  '
$ws.Range("D11").Value = '''0.07380'
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("D12").Value = '1.678.47'
$ws.Range("E12").Value = '  -5.20%  '
$ws.Range("D13").Value = '''4.555'
$ws.Range("E13").Value = '  -4.00%  '
$ws.Range("D14").Value = '''0.5770'
$ws.Range("E14").Value = '  -4.68%  '
$ws.Range("D15").Value = '1.908.41'
$ws.Range("E15").Value = '  -5.19%  '
$ws.Range("D16").Value = '''0.000008564'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = '''65.14'
$ws.Range("E17").Value = '  -12.53%  '
$ws.Range("D18").Value = '26.318.76'
$ws.Range("E18").Value = '  -7.43%  '
$ws.Range("D19").Value = '''5.008'
$ws.Range("E19").Value = '  -6.41%  '
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").Value = '''10.89'
$ws.Range("E21").Value = '  -4.12%  '
$ws.Range("D22").Value = '''187.05'
$ws.Range("E22").Value = '  -9.28%  '
$ws.Range("D23").Value = '''6.228'
$ws.Range("E23").Value = '  -7.30%  '
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("D25").Value = '''143.71'
$ws.Range("E25").Value = '  -4.96%  '
$ws.Range("D26").Value = '''7.577'
$ws.Range("E26").Value = '  -6.40%  '
$ws.Range("D27").Value = '''0.1179'
$ws.Range("E27").Value = '  -5.27%  '
$ws.Range("D28").Value = '''15.69'
$ws.Range("E28").Value = '  -3.29%  '
$ws.Range("E29").Value = '  -5.39%  '
$ws.Range("D30").Value = '''0.05811'
$ws.Range("E30").Value = '  -5.42%  '
$ws.Range("D31").Value = '''1.324'
$ws.Range("E31").Value = '  -6.05%  '
$ws.Range("E32").Value = '  -6.21%  '
$ws.Range("D33").Value = '''3.504'
$ws.Range("E33").Value = '  -6.44%  '
$ws.Range("D34").Value = '''1.666'
$ws.Range("E34").Value = '  -0.23%  '
$ws.Range("D35").Value = '''1.004'
$ws.Range("E35").Value = '  -3.44%  '
$ws.Range("D36").Value = '''0.5996'
$ws.Range("E36").Value = '  -5.46%  '
$ws.Range("D38").Value = '''2.658'
$ws.Range("E38").Value = '  -0.46%  '
$ws.Range("D39").Value = '1.098.58'
$ws.Range("E39").Value = '  -2.96%  '
$ws.Range("D40").Value = '''0.01607'
$ws.Range("E40").Value = '  -3.98%  '
$ws.Range("D41").Value = '''5.904'
$ws.Range("E41").Value = '  -6.36%  '
$ws.Range("D42").Value = '''0.8623'
$ws.Range("E42").Value = '  -0.96%  '
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").Value = '''99.49'
$ws.Range("D45").Value = '1.826.45'
$ws.Range("E45").Value = '  -5.17%  '
$ws.Range("D46").Value = '''0.00000000114'
$ws.Range("E46").Value = '  +4.04%  '
$ws.Range("D47").Value = '''56.40'
$ws.Range("E47").Value = '  -5.13%  '
$ws.Range("D48").Value = '''1.005'
$ws.Range("E48").Value = '  -0.24%  '
$ws.Range("D49").Value = '''8.057'
$ws.Range("E49").Value = '  -3.25%  '
$ws.Range("D50").Value = '''0.4315'
$ws.Range("E50").Value = '  -3.30%  '
$ws.Range("D51").Value = '''0.05214'
$ws.Range("E51").Value = '  -3.57%  '
